$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the raw input data (B2, C2, D2) ---
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 92
$ws.Range("D2").Value = 40000

# --- Center-align the whole table (A1:F7), including the two new blank
#     spacer rows (3 and 4) that get pulled into the used range ---
$ws.Range("A1:F7").HorizontalAlignment = -4108

# --- Conditional formatting: flag E7 in red when it drops below D2/1000 ---
$cf = $ws.Range("E7").FormatConditions.Add(1, 6, "=`$D`$2/1000")
$cf.Font.Color = 393372
$cf.Interior.Color = 13551615

# --- Move the active selection to B7 ---
$ws.Range("B7").Select() | Out-Null
